# Append a new arrivals row (row 20, record #19) to the "Main Data" sheet.
# Mirrors the existing rows: a Warsaw/LOT flight (2:30 PM scheduled) that
# landed as (SP-LIA) at 2:20 PM, 10 minutes early, on Friday, Jan 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = 19
$ws.Cells.Item($row, 2).Value = "Friday, Jan 13"
$ws.Cells.Item($row, 3).Value = "2:30 PM"
$ws.Cells.Item($row, 4).Value = "LO3993"
$ws.Cells.Item($row, 5).Value = "Warsaw"
$ws.Cells.Item($row, 6).Value = "(WAW)"
$ws.Cells.Item($row, 7).Value = "LOT "
$ws.Cells.Item($row, 8).Value = "E75S"
$ws.Cells.Item($row, 9).Value = "(SP-LIA)"
$ws.Cells.Item($row, 10).Value = "2:20 PM"

# Columns K and M stay blank in every other row of this table, but the
# underlying cell still exists in the sheet (it just has no value/type).
# Touching a no-op formatting property materialises the empty cell without
# creating a brand-new style entry, matching rows 2-19.
$ws.Cells.Item($row, 11).Borders.LineStyle = -4142
$ws.Cells.Item($row, 12).Value = "0 hours, -10 minutes"
$ws.Cells.Item($row, 13).Borders.LineStyle = -4142
